# Apply trade #11 (MarketMaking, closed 2026-02-17 13:08:38) to the workbook.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")

# ---- Summary sheet ----
$wsSummary.Range("B3").Value = 1199.62
$wsSummary.Range("B4").Value = -0.38
$wsSummary.Range("B5").Value = -0.69
$wsSummary.Range("B6").Value = 11
$wsSummary.Range("B7").Value = 4
$wsSummary.Range("B9").Value = 36.36

# ---- Strategy Status sheet (MarketMaking row) ----
$wsStrategy.Range("C4").Value = 99.62
$wsStrategy.Range("D4").Value = 11
$wsStrategy.Range("E4").Value = -0.38
$wsStrategy.Range("F4").Value = -0.38
$wsStrategy.Range("G4").Value = 36.36

# ---- Helper to append the new trade row on a trades-log sheet ----
function Add-TradeRow($ws) {
    $ws.Range("A12").Value = 11

    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = "2026-02-17"
    $ws.Range("B12").Style = "Normal"

    $ws.Range("C12").NumberFormat = "@"
    $ws.Range("C12").Value = "13:08:38"
    $ws.Range("C12").Style = "Normal"

    $ws.Range("D12").Value = "MarketMaking"
    $ws.Range("E12").Value = "DOWN"
    $ws.Range("F12").Value = 0.96
    $ws.Range("G12").Value = 0.97
    $ws.Range("H12").Value = "CLOSED"
    $ws.Range("I12").Value = 1.0417
    $ws.Range("J12").Value = 0.01
    $ws.Range("K12").Value = 99.62
    $ws.Range("L12").Value = 0
    $ws.Range("M12").Value = 0
    $ws.Range("N12").Value = 0.6
    $ws.Range("O12").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P12").Value = "early_exit"
    $ws.Range("Q12").Value = 0.13
}

Add-TradeRow $wsAllTrades
Add-TradeRow $wsMarketMaking
